# feat: added execute condition to test data
# fixes: ExcelReader to read the execute column
#
# This run regenerated the scraped "Results" / "CourseDetails" sheets:
#  - every sheet got a fresh timestamped name
#  - four of the "Results" sheets (positions 3,4,7,8) turned out to be
#    "CourseDetails" runs instead, so they are rewritten with the 4-row
#    Course 1 / Course 2 comparison data (and the execute/column bug is
#    fixed: Course 2's title + duration now reflect the real scrape)
#  - the sheet that used to hold CourseDetails data (position 12) was
#    actually a Results run, so it is rewritten with the full 48-row
#    Language/Level table
#  - three Results sheets (positions 5,6,9) had their B1 header stuck on
#    "Language" instead of "Level" (the execute-column read bug) - fixed

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename every sheet (position-based) to the new run names.
# ---------------------------------------------------------------------
$newNames = @(
    "edge_Results_20250903_094411",
    "chrome_Results_20250903_094415",
    "edge_CourseDetails_20250903_094",
    "chrome_CourseDetails_20250903_0",
    "edge_Results_20250903_100017",
    "chrome_Results_20250903_100034",
    "edge_CourseDetails_20250903_100",
    "chrome_CourseDetails_20250903_1",
    "chrome_Results_20250903_100355",
    "edge_Results_20250903_100409",
    "edge_Results_20250903_100724",
    "chrome_Results_20250903_100726"
)

for ($i = 1; $i -le $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

# ---------------------------------------------------------------------
# 2. Sheets at positions 3, 4, 7, 8 actually hold CourseDetails data:
#    shrink them from the 48-row Results layout down to the 4-row
#    Course 1 / Course 2 comparison.
# ---------------------------------------------------------------------
$courseDetailRows = @(
    @("Course 1", "Course 2"),
    @("title : Introduction to Front-End Development", "title : HTML, CSS, and Javascript for Web Developers"),
    @("rating : 4.8", "rating : 4.7"),
    @("duration : 1 - 4 Weeks", "duration : 1 - 3 Months")
)

foreach ($pos in 3, 4, 7, 8) {
    $ws = $wb.Worksheets.Item($pos)

    # drop the old rows 5-48 so the sheet shrinks back to A1:B4
    $ws.Range("A5:B48").EntireRow.Delete()

    for ($r = 1; $r -le $courseDetailRows.Length; $r++) {
        $pair = $courseDetailRows[$r - 1]
        $ws.Cells.Item($r, 1).Value = $pair[0]
        $ws.Cells.Item($r, 2).Value = $pair[1]
    }
}

# ---------------------------------------------------------------------
# 3. The sheet at position 12 actually holds Results data: expand it
#    from the 4-row CourseDetails layout into the full 48-row
#    Language/Level table.
# ---------------------------------------------------------------------
$resultsPairedRows = @(
    @("Language", "Level"),
    @("English(215)", "Beginner(98)"),
    @("Chinese(186)", "Intermediate(91)"),
    @("Spanish(177)", "Advanced(8)"),
    @("French(170)", "Mixed(50)")
)

$resultsSingleRows = @(
    "Arabic(166)", "German(166)", "Portuguese(166)", "Russian(166)",
    "Hindi(165)", "Indonesian(165)", "Korean(154)", "Japanese(153)",
    "Italian(152)", "Swedish(152)", "Dutch(151)", "Greek(151)",
    "Kazakh(151)", "Polish(151)", "Thai(151)", "Turkish(151)",
    "Ukrainian(151)", "Pushto(136)", "Vietnamese(122)", "Urdu(113)",
    "Azerbaijani(109)", "Bengali(109)", "Hungarian(109)", "Persian(23)",
    "Afrikaans(16)", "Oriya(9)", "Mongolian(7)", "Catalan(6)",
    "Slovak(5)", "Tamil(5)", "Serbian(4)", "Uzbek(4)",
    "Czech(3)", "Hebrew(3)", "Lithuanian(3)", "Romanian(3)",
    "Albanian(2)", "Burmese(2)", "Javanese(2)", "Bulgarian(1)",
    "Georgian(1)", "Marathi(1)", "Swahili(1)"
)

$wsResults = $wb.Worksheets.Item(12)

for ($r = 1; $r -le $resultsPairedRows.Length; $r++) {
    $pair = $resultsPairedRows[$r - 1]
    $wsResults.Cells.Item($r, 1).Value = $pair[0]
    $wsResults.Cells.Item($r, 2).Value = $pair[1]
}

$row = $resultsPairedRows.Length + 1
foreach ($lang in $resultsSingleRows) {
    $wsResults.Cells.Item($row, 1).Value = $lang
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 4. Sheets at positions 5, 6, 9 had the execute-column read bug: B1
#    held "Language" (copied from A1) instead of "Level". Fix it.
# ---------------------------------------------------------------------
foreach ($pos in 5, 6, 9) {
    $ws = $wb.Worksheets.Item($pos)
    $ws.Range("B1").Value = "Level"
}
